$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("Type") to hold the new "Variable" column.
$ws.Columns.Item(2).Insert()

# New header + value for the inserted "Variable" column.
$ws.Range("B1").Value = "Variable"
$ws.Range("B2").Value = "e1"

# Give the new column a sensible width (close to the author's saved width).
$ws.Columns.Item(2).ColumnWidth = 18.721354166666668

# Match the saved selection state.
$ws.Range("B3").Select()
